# This workbook's sheet is protected (sheetProtection with a legacy
# password hash). Unprotecting/reprotecting via Unprotect()/Protect()
# would strip that legacy hash (the engine always re-hashes with a
# modern SHA-512 scheme), which would be an unwanted, unrelated change.
# Instead, momentarily flip Locked off on just the cells we need to
# touch, write the new values, then flip Locked back on -- this keeps
# the worksheet's <sheetProtection> element completely untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update the confidential disclosure date notice (row 10, col A) ---
$ws.Range("A10").Locked = $false
$ws.Range("A10").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-07-09 for illustrative purposes only and are subject to change."
# Re-assigning a multi-line value makes the engine auto-grow row 10's
# height; AutoFit puts it straight back to the sheet's implicit default
# so no spurious ht/customHeight attributes are written out.
$ws.Rows.Item(10).AutoFit()
$ws.Range("A10").Locked = $true

# --- Update Weight (D) / Percent Change (E) for rows 2-7 ---
$ws.Range("D2:E7").Locked = $false

$ws.Range("D2").Value = 0.2643883625788284
$ws.Range("E2").Value = 0.007718431614695875

$ws.Range("D3").Value = 0.5296920428812297
$ws.Range("E3").Value = 0.01496126102057183

$ws.Range("D4").Value = 0.05270232356019872
$ws.Range("E4").Value = 0.01101321585903081

$ws.Range("D5").Value = 0.09570739918805964
$ws.Range("E5").Value = 0.02215988779803646

$ws.Range("D6").Value = 0.05750987179168351
$ws.Range("E6").Value = 0.02069122328331074

$ws.Range("D7").Value = 0.9999999999999999
$ws.Range("E7").Value = 0.01385676130167846

$ws.Range("D2:E7").Locked = $true
